$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "79.745.14"
$ws.Range("D3").Value = "3.198.36"
$ws.Range("E3").Value = "  +4.91%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'205.22"
$ws.Range("E5").Value = "  +1.72%  "
$ws.Range("D6").Value = "'635.05"
$ws.Range("E6").Value = "  +1.46%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.237"
$ws.Range("E8").Value = "  +13.18%  "
$ws.Range("E9").Value = "  +5.57%  "
$ws.Range("D10").Value = "3.197.63"
$ws.Range("E10").Value = "  +4.92%  "
$ws.Range("D11").Value = "'0.579"
$ws.Range("E11").Value = "  +32.25%  "
$ws.Range("D12").Value = "'0.165"
$ws.Range("E12").Value = "  +2.99%  "
$ws.Range("D13").Value = "'5.51"
$ws.Range("E13").Value = "  +6.97%  "
$ws.Range("D14").Value = "'0.0000232"
$ws.Range("E14").Value = "  +19.11%  "
$ws.Range("D15").Value = "3.785.69"
$ws.Range("E15").Value = "  +4.75%  "
$ws.Range("D16").Value = "'31.86"
$ws.Range("E16").Value = "  +8.50%  "
$ws.Range("D17").Value = "79.552.15"
$ws.Range("E17").Value = "  +4.15%  "
$ws.Range("D18").Value = "3.203.04"
$ws.Range("E18").Value = "  +4.90%  "
$ws.Range("D19").Value = "'14.48"
$ws.Range("E19").Value = "  +6.76%  "
$ws.Range("D20").Value = "'3.02"
$ws.Range("E20").Value = "  +30.12%  "
$ws.Range("D21").Value = "'9.19"
$ws.Range("E21").Value = "  +0.65%  "
$ws.Range("D22").Value = "'429.42"
$ws.Range("E22").Value = "  +14.43%  "
$ws.Range("D23").Value = "'5.12"
$ws.Range("E23").Value = "  +17.23%  "
$ws.Range("D24").Value = "'11.26"
$ws.Range("E24").Value = "  +13.37%  "
$ws.Range("D25").Value = "3.364.47"
$ws.Range("E25").Value = "  +4.89%  "
$ws.Range("D26").Value = "'77.09"
$ws.Range("E26").Value = "  +4.65%  "
$ws.Range("D27").Value = "'4.73"
$ws.Range("E27").Value = "  +6.75%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").Value = "'0.0000119"
$ws.Range("E29").Value = "  +6.54%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").Value = "'9.01"
$ws.Range("E31").Value = "  +8.17%  "
$ws.Range("E32").Value = "  +5.08%  "
$ws.Range("D33").Value = "'527.13"
$ws.Range("E33").Value = "  +4.06%  "
$ws.Range("E34").Value = "  +1.55%  "
$ws.Range("E35").Value = "  +25.39%  "
$ws.Range("D36").Value = "'23.00"
$ws.Range("E36").Value = "  +10.42%  "
$ws.Range("D37").Value = "'0.120"
$ws.Range("E37").Value = "  +12.70%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").Value = "'0.405"
$ws.Range("E39").Value = "  +4.96%  "
$ws.Range("D40").Value = "'165.35"
$ws.Range("E40").Value = "  +1.37%  "
$ws.Range("D41").Value = "'20.03"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "'192.36"
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").Value = "'5.52"
$ws.Range("E44").Value = "  +6.12%  "
$ws.Range("D45").Value = "'0.822"
$ws.Range("E45").Value = "  +3.01%  "
$ws.Range("E46").Value = "  +7.29%  "
$ws.Range("E47").Value = "  +3.59%  "
$ws.Range("D48").Value = "'43.30"
$ws.Range("E48").Value = "  +2.33%  "
$ws.Range("D49").Value = "'25.80"
$ws.Range("E49").Value = "  +14.88%  "
$ws.Range("E50").Value = "  +4.65%  "
$ws.Range("D51").Value = "'2.52"
$ws.Range("E51").Value = "  +1.61%  "
